# Ajustes do projeto para executar local.
# Updates the "Escopo" rubric sheet: marks several previously
# "Em Andamento" items as "Concluido", fills in blank status cells,
# adds a "Não feito" status, and bumps the related percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Escopo")

# --- Row 6: Designer de aplicação utilizando DDD -> now Concluido / 100% ---
$ws.Range("E6").Value = "Concluido"
$ws.Range("F6").Value = 1

# --- Rows 9, 10, 11, 12: flip "Em Andamento" -> "Concluido" ---
$ws.Range("E9").Value = "Concluido"
$ws.Range("E10").Value = "Concluido"
$ws.Range("E11").Value = "Concluido"
$ws.Range("E12").Value = "Concluido"

# --- Row 13: flip status + percentage to complete ---
$ws.Range("E13").Value = "Concluido"
$ws.Range("F13").Value = 1

# --- Rows 14, 15: fill in blank Status with Concluido ---
$ws.Range("E14").Value = "Concluido"
$ws.Range("E15").Value = "Concluido"

# --- Row 16: fill in blank Status with the new "Não feito" value ---
$ws.Range("E16").Value = "Não feito"

# --- Rows 17 - 22: fill in blank Status with Concluido ---
$ws.Range("E17").Value = "Concluido"
$ws.Range("E18").Value = "Concluido"
$ws.Range("E19").Value = "Concluido"
$ws.Range("E20").Value = "Concluido"
$ws.Range("E21").Value = "Concluido"
$ws.Range("E22").Value = "Concluido"

# --- Row 24: fill in Status + bump percentage to complete ---
$ws.Range("E24").Value = "Concluido"
$ws.Range("F24").Value = 1

# --- Row 25: flip "Em Andamento" -> "Concluido" ---
$ws.Range("E25").Value = "Concluido"

# Update the active selection to match the author's last position.
$ws.Range("E18").Select() | Out-Null
